$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B to make room for the two new
# "as-of" dates (Jun_26 appears twice, Jun_27 once). The existing
# Jun_17/Jun_15/Jun_13/Jun_10 columns (B:E) shift right to E:H.
$ws.Columns("B:D").Insert()

# Keep the "8.0" custom width that every date column (C:H) already used
# before the insert (Excel expresses this in character units, which is
# ~0.8333 narrower than the stored worksheet width at the default font).
$ws.Columns("C:H").ColumnWidth = 7.166666666666667

# New header row values for the freshly inserted date columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the new columns (B:G) with "UN" for every analyst data row (2-27).
# Column H keeps whatever rating/price-target text already lived there
# (it was column E before the insert, so it moved automatically).
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("B$r").Value = "UN"
    $ws.Range("C$r").Value = "UN"
    $ws.Range("D$r").Value = "UN"
    $ws.Range("E$r").Value = "UN"
    $ws.Range("F$r").Value = "UN"
    $ws.Range("G$r").Value = "UN"
}

# Add a new analyst group at the bottom of the sheet.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
